$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("C3").Value = "CityCode,LandOfficeCode"
$ws.Range("E9").Value = 4

$ws.Activate()
$ws.Range("G3").Select()
